# Asset & Rental : implement space
# Replace the "office" placeholder value in P2/P3 (location column) with
# actual row numbers "1" and "2" respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2").Value = "1"
$ws.Range("P3").Value = "2"

$ws.Range("P3").Select()
